$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "MSG: None

MSG: The decision to acquire the rights for `"Barbie`" has been recorded.
"
$ws.Range("D2").Value = "Barbie_was_selected, "
$ws.Range("C3").Value = "MSG: None

MSG: The rights to both movies have been acquired.
"
$ws.Range("D3").Value = "both_movies, "
$ws.Range("C4").Value = "MSG: None

MSG: The decision has been recorded to acquire rights for `"Barbie.`"
"
$ws.Range("D4").Value = "Barbie_was_selected, "
$ws.Range("C5").Value = "MSG: None

MSG: The decision has been made to acquire the rights for `"Barbie.`"
"
$ws.Range("D5").Value = "Barbie_was_selected, "
$ws.Range("C6").Value = "MSG: None

MSG: The rights for both movies have been successfully acquired.
"
$ws.Range("D6").Value = "both_movies, "
$ws.Range("C7").Value = "MSG: None

MSG: The decision-making process did not result in a clear choice for Friday's movie, thus we ended up with no decision.
"
$ws.Range("D7").Value = "no_decision, "
$ws.Range("C8").Value = "MSG: None

MSG: The decision has been recorded that no movie was selected for Friday.
"
$ws.Range("D8").Value = "no_decision, "
$ws.Range("C9").Value = "MSG: None

MSG: I have successfully recorded the decision to acquire the rights for both movies.
"
$ws.Range("D9").Value = "both_movies, "
$ws.Range("C10").Value = "MSG: None

MSG: The decision has been recorded as no decision was reached regarding the movie to be shown on Friday.
"
$ws.Range("D10").Value = "no_decision, "
$ws.Range("C11").Value = "MSG: None

MSG: The decision-making process concluded without an agreement on which movie to show on Friday, resulting in no decision about selecting a movie.
"
$ws.Range("D11").Value = "no_decision, "
$ws.Range("C12").Value = "MSG: None

MSG: The decision about which movie to show on Friday was not made.
"
$ws.Range("D12").Value = "no_decision, "
$ws.Range("C13").Value = "MSG: None

MSG: The decision to acquire the rights for `"Oppenheimer`" has been successfully recorded.
"
$ws.Range("D13").Value = "Oppenheimer_was_selected, "
$ws.Range("C14").Value = "MSG: None

MSG: The decision has been recorded as no decision being made about which movie to show on Friday.
"
$ws.Range("D14").Value = "no_decision, "
$ws.Range("C15").Value = "MSG: None

MSG: The decision has been made to acquire the rights for `"Barbie.`"
"
$ws.Range("D15").Value = "Barbie_was_selected, "
$ws.Range("C16").Value = "MSG: None

MSG: The decision has been recorded, and no movie will be acquired for Friday.
"
$ws.Range("D16").Value = "no_decision, "
$ws.Range("C17").Value = "MSG: None

MSG: The decision for Friday's movie has resulted in no consensus, so no movie will be acquired at this time.
"
$ws.Range("D17").Value = "no_decision, "
$ws.Range("C18").Value = "MSG: None

MSG: The decision has been made to acquire rights for both movies.
"
$ws.Range("D18").Value = "both_movies, "
$ws.Range("C19").Value = "MSG: None

MSG: The decision has been made to acquire `"Barbie`" for the movie to be shown on Friday.
"
$ws.Range("D19").Value = "Barbie_was_selected, "
$ws.Range("C20").Value = "MSG: None

MSG: The decision to acquire the rights for the movie `"Barbie`" has been successfully recorded.
"
$ws.Range("D20").Value = "Barbie_was_selected, "
$ws.Range("C21").Value = "MSG: None

MSG: The decision about the movie to be shown on Friday could not be made.
"
$ws.Range("D21").Value = "no_decision, "
$ws.Range("C22").Value = "MSG: None

MSG: The decision regarding the movie for Friday has resulted in no conclusion.
"
$ws.Range("D22").Value = "no_decision, "
$ws.Range("C23").Value = "MSG: None

MSG: The decision to show both movies has been recorded successfully.
"
$ws.Range("D23").Value = "both_movies, "
$ws.Range("C24").Value = "MSG: None

MSG: The decision has been made to acquire the rights for `"Barbie.`"
"
$ws.Range("D24").Value = "Barbie_was_selected, "
$ws.Range("C25").Value = "MSG: None

MSG: The rights for both movies have been successfully acquired.
"
$ws.Range("D25").Value = "both_movies, "
$ws.Range("C26").Value = "MSG: None

MSG: The decision has been made to acquire the rights for the movie `"Barbie.`"
"
$ws.Range("D26").Value = "Barbie_was_selected, "
$ws.Range("C27").Value = "MSG: None

MSG: The committee did not reach a decision about which movie to show on Friday.
"
$ws.Range("D27").Value = "no_decision, "
$ws.Range("C28").Value = "MSG: None

MSG: The decision has been made to acquire the rights for `"Barbie.`"
"
$ws.Range("D28").Value = "Barbie_was_selected, "
$ws.Range("C29").Value = "MSG: None

MSG: The decision has been recorded, and the rights to `"Barbie`" will be acquired.
"
$ws.Range("D29").Value = "Barbie_was_selected, "
$ws.Range("C30").Value = "MSG: None

MSG: The decision to acquire the rights for both movies has been recorded.
"
$ws.Range("D30").Value = "both_movies, "
$ws.Range("C31").Value = "MSG: None

MSG: The decision regarding the movie to show on Friday resulted in no agreement.
"
$ws.Range("D31").Value = "no_decision, "
$ws.Range("C32").Value = "MSG: None

MSG: The decision to show both `"Oppenheimer`" and `"Barbie`" on Friday has been made, and the rights to both movies will be acquired.
"
$ws.Range("D32").Value = "both_movies, "
$ws.Range("C33").Value = "MSG: None

MSG: The rights to `"Barbie`" have been acquired for the movie to be shown on Friday.
"
$ws.Range("D33").Value = "Barbie_was_selected, "
$ws.Range("C34").Value = "MSG: None

MSG: The decision has been made to select `"Barbie`" as the movie for Friday.
"
$ws.Range("D34").Value = "Barbie_was_selected, "
$ws.Range("C35").Value = "MSG: None

MSG: I have successfully recorded the decision to acquire rights for both movies.
"
$ws.Range("D35").Value = "both_movies, "
$ws.Range("C36").Value = "MSG: None

MSG: The decision has been recorded, and no movie has been selected for Friday.
"
$ws.Range("D36").Value = "no_decision, "
$ws.Range("C37").Value = "MSG: None

MSG: I have recorded the decision to acquire the rights for `"Oppenheimer.`"
"
$ws.Range("D37").Value = "Oppenheimer_was_selected, "
$ws.Range("C38").Value = "MSG: None

MSG: The rights to both movies have been acquired for Friday's showing.
"
$ws.Range("D38").Value = "both_movies, "
$ws.Range("C39").Value = "MSG: None

MSG: The decision has been recorded to acquire rights for `"Oppenheimer`" to be shown on Friday.
"
$ws.Range("D39").Value = "Oppenheimer_was_selected, "
$ws.Range("C40").Value = "MSG: None

MSG: The decision to show a movie on Friday was not clearly finalized, so no action will be taken regarding the acquisition of movie rights.
"
$ws.Range("D40").Value = "no_decision, "
$ws.Range("C41").Value = "MSG: None

MSG: The decision has been recorded, and the rights for `"Barbie`" have been acquired.
"
$ws.Range("D41").Value = "Barbie_was_selected, "
$ws.Range("C42").Value = "MSG: None

MSG: The decision regarding the movie for Friday has ultimately resulted in no decision being made.
"
$ws.Range("D42").Value = "no_decision, "
$ws.Range("C43").Value = "MSG: None

MSG: The decision regarding the movie to be shown on Friday has resulted in no consensus. Therefore, no movie will be acquired at this time.
"
$ws.Range("D43").Value = "no_decision, "
$ws.Range("C44").Value = "MSG: None

MSG: The movie rights to both `"Oppenheimer`" and `"Barbie`" have been successfully acquired.
"
$ws.Range("D44").Value = "both_movies, "
$ws.Range("C45").Value = "MSG: None

MSG: The decision has been recorded as no_decision.
"
$ws.Range("D45").Value = "no_decision, "
$ws.Range("C46").Value = "MSG: None

MSG: The decision to acquire the rights for `"Barbie`" has been successfully recorded.
"
$ws.Range("D46").Value = "Barbie_was_selected, "
$ws.Range("C47").Value = "MSG: None

MSG: The decision to show `"Oppenheimer`" has been recorded.
"
$ws.Range("D47").Value = "Oppenheimer_was_selected, "
$ws.Range("C48").Value = "MSG: None

MSG: The decision has been recorded, and `"Barbie`" will be the movie shown on Friday.
"
$ws.Range("D48").Value = "Barbie_was_selected, "
$ws.Range("C49").Value = "MSG: None

MSG: The decision has been recorded as no agreement was reached on a movie to be shown on Friday.
"
$ws.Range("D49").Value = "no_decision, "
$ws.Range("C50").Value = "MSG: None

MSG: The committee has not reached a decision about which movie to show on Friday, as they are planning to vote to determine the choice. Therefore, I will call the no_decision function.
"
$ws.Range("D50").Value = "no_decision, "
$ws.Range("C51").Value = "MSG: None

MSG: The decision has been recorded as no decision regarding Friday's movie can be made.
"
$ws.Range("D51").Value = "no_decision, "
$ws.Range("C52").Value = "MSG: None

MSG: The decision has been made, and no movie will be shown on Friday.
"
$ws.Range("D52").Value = "no_decision, "
$ws.Range("C53").Value = "MSG: None

MSG: The decision to acquire the rights for both movies has been successfully recorded.
"
$ws.Range("D53").Value = "both_movies, "
$ws.Range("C54").Value = "MSG: None

MSG: The decision has been recorded to acquire the rights for both movies.
"
$ws.Range("D54").Value = "both_movies, "
$ws.Range("C55").Value = "MSG: None

MSG: The decision has been recorded as no final agreement was reached regarding the movie for Friday.
"
$ws.Range("D55").Value = "no_decision, "
$ws.Range("C56").Value = "MSG: None

MSG: The decision about which movie to show on Friday has not been made.
"
$ws.Range("D56").Value = "no_decision, "
$ws.Range("C57").Value = "MSG: None

MSG: I have acquired the rights to both movies as that was the committee's final decision.
"
$ws.Range("D57").Value = "both_movies, "
$ws.Range("C58").Value = "MSG: None

MSG: The function has been called, and no decision has been made regarding the movie for Friday.
"
$ws.Range("D58").Value = "no_decision, "
$ws.Range("C59").Value = "MSG: None

MSG: The decision has been successfully recorded, and `"Barbie`" has been selected to be shown on Friday.
"
$ws.Range("D59").Value = "Barbie_was_selected, "
$ws.Range("C60").Value = "MSG: None

MSG: The decision for Friday's movie resulted in no agreement.
"
$ws.Range("D60").Value = "no_decision, "
$ws.Range("C61").Value = "MSG: None

MSG: The committee did not come to a decision regarding the movie to show on Friday.
"
$ws.Range("D61").Value = "no_decision, "
$ws.Range("C62").Value = "MSG: None

MSG: The decision has been recorded, and no movie will be acquired for Friday.
"
$ws.Range("D62").Value = "no_decision, "
$ws.Range("C63").Value = "MSG: None

MSG: The rights for both movies have been acquired successfully.
"
$ws.Range("D63").Value = "both_movies, "
$ws.Range("C64").Value = "MSG: None

MSG: The rights to both movies have been acquired successfully.
"
$ws.Range("D64").Value = "both_movies, "
$ws.Range("C65").Value = "MSG: None

MSG: The decision-making process did not lead to a selection for Friday's movie.
"
$ws.Range("D65").Value = "no_decision, "
$ws.Range("C66").Value = "MSG: None

MSG: The decision has been recorded to acquire the rights for `"Oppenheimer.`"
"
$ws.Range("D66").Value = "Oppenheimer_was_selected, "
$ws.Range("C67").Value = "MSG: None

MSG: The decision has been recorded, indicating that no movie was selected for Friday.
"
$ws.Range("D67").Value = "no_decision, "
$ws.Range("C68").Value = "MSG: None

MSG: The decision has been recorded. `"Barbie`" will be the movie shown on Friday.
"
$ws.Range("D68").Value = "Barbie_was_selected, "
$ws.Range("C69").Value = "MSG: None

MSG: The decision has been recorded, and no movie was selected for Friday.
"
$ws.Range("D69").Value = "no_decision, "
